$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BasePath_BP_File")
$ws.Activate()

# Add the two new BP rows (Staff / Customer AI Character) right after the
# existing data, mirroring the existing table layout (Id / Id / FString name).
$ws.Range("A86").Value2 = 9001
$ws.Range("B86").Value2 = 9001
$ws.Range("C86").Value2 = "BP_StaffAICharacter"

$ws.Range("A87").Value2 = 9002
$ws.Range("B87").Value2 = 9001
$ws.Range("C87").Value2 = "BP_CustomerAICharacter"

# The name column (C) uses the same cell style as the header's "FString"
# type cell (C2) -- copy that formatting onto the two new name cells.
$ws.Range("C2").Copy()
$ws.Range("C86:C87").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Match the author's final selection/scroll position.
[void]$ws.Range("A88").Select()
